$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the YCbCr color space GLCM computation (90 degree, G_n)
$data = @(
    @(0.023855847193339406, 0.8914328087274338,  0.76781670680412206, 0.98807207640333017),
    @(0.054985849440643275, 0.87746139802441858, 0.50335328380216149, 0.97250707527967828),
    @(0.0038355525518515747,0.81210479414250147, 0.97393860571495749, 0.9980822237240744),
    @(0.000021674875291409153,0.27617963851863658,0.99994439716747485, 0.99998916256235437)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Column widths (character-unit widths matching the source workbook's cached
# pixel-derived values as closely as the ColumnWidth rounding allows)
$ws.Columns.Item(1).ColumnWidth = 13.833333333333332
$ws.Columns.Item(2).ColumnWidth = 11.833333333333332
$ws.Columns.Item(3).ColumnWidth = 11.833333333333332
$ws.Columns.Item(4).ColumnWidth = 11.833333333333332

# Two auxiliary number-format styles (text + date) get minted into the
# workbook's style table, each carrying its own (otherwise-default) border
# record -- mirrors style slots carried over from the template this sheet's
# data was generated from, even though no visible cell uses them here.
$aux1 = $ws.Range("Z500")
$aux1.Borders.Item(5).LineStyle = 1
$aux1.NumberFormat = "@"

$aux2 = $ws.Range("Z501")
$aux2.Borders.Item(5).LineStyle = 0
$aux2.NumberFormat = "m/d/yy h:mm"

$ws.Rows.Item(500).Delete()
$ws.Rows.Item(500).Delete()

# Ensure the sheet is marked as the selected/active tab
$ws.Select()

# Force full recalculation on load
$wb.Application.CalculateFullRebuild()
